# Add "tei_quote" block-quote styling, consistently basing it on the
# built-in "Normal (Web)" style (and its linked character style) so the
# TEI-to-docx conversion renders block quotes the same way everywhere.

$d = $word.ActiveDocument

# --- tei_quote --------------------------------------------------------
$teiQuote = $d.Styles.Add("teiquote", 1)
$teiQuote.NameLocal = "tei_quote"
$teiQuote.BaseStyle = "NormalWeb"
$teiQuote.LinkStyle = "teiquoteChar"
$teiQuote.QuickStyle = $true
$teiQuote.ParagraphFormat.SpaceBefore = 6
$teiQuote.ParagraphFormat.SpaceBeforeAuto = $false
$teiQuote.ParagraphFormat.SpaceAfter = 6
$teiQuote.ParagraphFormat.LineSpacingRule = 0
$teiQuote.ParagraphFormat.LineSpacing = 12
$teiQuote.ParagraphFormat.LeftIndent = 8.5
$teiQuote.ParagraphFormat.RightIndent = 8.5
$teiQuote.Font.Size = 11

# --- Normal (Web) -------------------------------------------------------
$normalWeb = $d.Styles.Add("Normal (Web)", 1)
$normalWeb.BaseStyle = "Normal"
$normalWeb.LinkStyle = "NormalWebChar"
$normalWeb.Priority = 99
$normalWeb.Font.Name = "Times New Roman"
$normalWeb.Font.Size = 12
$normalWeb.Font.SizeBi = 12

# --- Normal (Web) Char ---------------------------------------------------
$normalWebChar = $d.Styles.Add("Normal (Web) Char", 2)
$normalWebChar.BaseStyle = "DefaultParagraphFont"
$normalWebChar.LinkStyle = "NormalWeb"
$normalWebChar.Priority = 99
$normalWebChar.Font.Name = "Times New Roman"
$normalWebChar.Font.Size = 12
$normalWebChar.Font.SizeBi = 12

# --- tei_quote Char -------------------------------------------------------
$teiQuoteChar = $d.Styles.Add("teiquoteChar", 2)
$teiQuoteChar.NameLocal = "tei_quote Char"
$teiQuoteChar.BaseStyle = "NormalWebChar"
$teiQuoteChar.LinkStyle = "teiquote"
$teiQuoteChar.Font.Name = "Times New Roman"
$teiQuoteChar.Font.Size = 12
$teiQuoteChar.Font.SizeBi = 12

Write-Output "tei_quote styles added"
